$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "302.19"
Set-TextValue $ws.Range("E2") "1.19%"
Set-TextValue $ws.Range("G2") "12"
Set-TextValue $ws.Range("E3") "1.22%"
Set-TextValue $ws.Range("G3") "12"
Set-TextValue $ws.Range("D4") "5.101"
Set-TextValue $ws.Range("E4") "0.09%"
Set-TextValue $ws.Range("G4") "12"
Set-TextValue $ws.Range("D5") "0.07839"
Set-TextValue $ws.Range("E5") "-2.27%"
Set-TextValue $ws.Range("G5") "12"
Set-TextValue $ws.Range("D6") "2.283"
Set-TextValue $ws.Range("E6") "-12.23%"
Set-TextValue $ws.Range("G6") "12"
Set-TextValue $ws.Range("D7") "7.807"
Set-TextValue $ws.Range("E7") "-0.14%"
Set-TextValue $ws.Range("G7") "12"
Set-TextValue $ws.Range("D8") "3.814"
Set-TextValue $ws.Range("E8") "-0.21%"
Set-TextValue $ws.Range("G8") "12"
Set-TextValue $ws.Range("D9") "0.9251"
Set-TextValue $ws.Range("E9") "0.73%"
Set-TextValue $ws.Range("G9") "12"
Set-TextValue $ws.Range("D10") "0.1774"
Set-TextValue $ws.Range("E10") "2.46%"
Set-TextValue $ws.Range("G10") "12"
Set-TextValue $ws.Range("D11") "0.07517"
Set-TextValue $ws.Range("E11") "2.19%"
Set-TextValue $ws.Range("G11") "12"
Set-TextValue $ws.Range("D12") "0.08945"
Set-TextValue $ws.Range("E12") "6.50%"
Set-TextValue $ws.Range("G12") "12"
Set-TextValue $ws.Range("D13") "0.03038"
Set-TextValue $ws.Range("E13") "0.55%"
Set-TextValue $ws.Range("G13") "12"
Set-TextValue $ws.Range("E14") "0.75%"
Set-TextValue $ws.Range("G14") "12"
Set-TextValue $ws.Range("D15") "0.001509"
Set-TextValue $ws.Range("E15") "0.59%"
Set-TextValue $ws.Range("G15") "12"
Set-TextValue $ws.Range("D16") "0.005789"
Set-TextValue $ws.Range("E16") "-2.95%"
Set-TextValue $ws.Range("G16") "12"
Set-TextValue $ws.Range("D17") "3.469"
Set-TextValue $ws.Range("E17") "-0.95%"
Set-TextValue $ws.Range("G17") "12"
Set-TextValue $ws.Range("D18") "2.252"
Set-TextValue $ws.Range("E18") "0.00%"
Set-TextValue $ws.Range("G18") "12"
Set-TextValue $ws.Range("G19") "12"
Set-TextValue $ws.Range("D20") "0.1335"
Set-TextValue $ws.Range("E20") "-0.15%"
Set-TextValue $ws.Range("G20") "12"
Set-TextValue $ws.Range("D21") "4.341"
Set-TextValue $ws.Range("E21") "-5.65%"
Set-TextValue $ws.Range("G21") "12"
Set-TextValue $ws.Range("D22") "0.1814"
Set-TextValue $ws.Range("E22") "13.47%"
Set-TextValue $ws.Range("G22") "12"
Set-TextValue $ws.Range("D23") "0.04582"
Set-TextValue $ws.Range("E23") "-0.53%"
Set-TextValue $ws.Range("G23") "12"
Set-TextValue $ws.Range("D24") "0.001247"
Set-TextValue $ws.Range("E24") "0.61%"
Set-TextValue $ws.Range("G24") "12"
Set-TextValue $ws.Range("D25") "0.004482"
Set-TextValue $ws.Range("E25") "0.73%"
Set-TextValue $ws.Range("G25") "12"
Set-TextValue $ws.Range("D26") "0.0001246"
Set-TextValue $ws.Range("E26") "4.69%"
Set-TextValue $ws.Range("G26") "12"
Set-TextValue $ws.Range("E27") "-1.31%"
Set-TextValue $ws.Range("G27") "12"
Set-TextValue $ws.Range("G28") "12"
Set-TextValue $ws.Range("G29") "12"
Set-TextValue $ws.Range("G30") "12"
Set-TextValue $ws.Range("G31") "12"
Set-TextValue $ws.Range("G32") "12"
Set-TextValue $ws.Range("G33") "12"
Set-TextValue $ws.Range("G34") "12"
Set-TextValue $ws.Range("G35") "12"
Set-TextValue $ws.Range("G36") "12"
Set-TextValue $ws.Range("G37") "12"
Set-TextValue $ws.Range("G38") "12"
Set-TextValue $ws.Range("D39") "0.01777"
Set-TextValue $ws.Range("E39") "-2.70%"
Set-TextValue $ws.Range("G39") "12"
Set-TextValue $ws.Range("D40") "0.04796"
Set-TextValue $ws.Range("E40") "5.27%"
Set-TextValue $ws.Range("G40") "12"
Set-TextValue $ws.Range("D41") "0.007351"
Set-TextValue $ws.Range("E41") "4.95%"
Set-TextValue $ws.Range("G41") "12"
Set-TextValue $ws.Range("D42") "0.1364"
Set-TextValue $ws.Range("E42") "1.62%"
Set-TextValue $ws.Range("G42") "12"
Set-TextValue $ws.Range("D43") "0.002118"
Set-TextValue $ws.Range("E43") "-5.43%"
Set-TextValue $ws.Range("G43") "12"
Set-TextValue $ws.Range("D44") "0.01051"
Set-TextValue $ws.Range("E44") "8.48%"
Set-TextValue $ws.Range("G44") "12"
Set-TextValue $ws.Range("D45") "0.00006281"
Set-TextValue $ws.Range("E45") "-4.09%"
Set-TextValue $ws.Range("G45") "12"
Set-TextValue $ws.Range("D46") "0.00000000748"
Set-TextValue $ws.Range("E46") "-0.31%"
Set-TextValue $ws.Range("G46") "12"
Set-TextValue $ws.Range("G47") "12"
Set-TextValue $ws.Range("D48") "0.7213"
Set-TextValue $ws.Range("E48") "-12.10%"
Set-TextValue $ws.Range("G48") "12"
Set-TextValue $ws.Range("D49") "0.00002093"
Set-TextValue $ws.Range("E49") "-0.31%"
Set-TextValue $ws.Range("G49") "12"
Set-TextValue $ws.Range("D50") "0.0001994"
Set-TextValue $ws.Range("E50") "-0.31%"
Set-TextValue $ws.Range("G50") "12"
Set-TextValue $ws.Range("G51") "12"
